# "Generate Report for Handoff"
#
# The localization-status report was regenerated: the handoff just kicked
# off, so the per-language status flips from the old "handed back" message
# to "Ready for handoff", and the two timestamps that track the latest
# generated/handoff xliff move forward a few seconds. Because the new
# status text is much shorter than the old one, the Status-ish columns
# that were sized to fit it shrink accordingly.

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item(1)   # "Overview"
$ws_zhcn     = $wb.Worksheets.Item(2)   # "zh-cn"
$ws_dede     = $wb.Worksheets.Item(3)   # "de-de"

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff" ---
$ws_overview.Range("E2").Value = "Ready for handoff"
$ws_overview.Range("F2").Value = "Ready for handoff"
$ws_zhcn.Range("C2").Value = "Ready for handoff"
$ws_dede.Range("C2").Value = "Ready for handoff"

# --- Timestamps ---
# Overview!G2 ("Latest HO Xliff Generate Date") and de-de!H2 ("Latest
# Handoff Datetime") shared the old value "2016-09-05 05:04:52" and both
# move to the new, later value "2016-09-05 05:05:44".
$ws_overview.Range("G2").Value = "2016-09-05 05:05:44"
$ws_dede.Range("H2").Value = "2016-09-05 05:05:44"

# zh-cn!H2 ("Latest Handoff Datetime") moves from "2016-09-05 05:04:48" to
# "2016-09-05 05:05:39".
$ws_zhcn.Range("H2").Value = "2016-09-05 05:05:39"

# --- Column widths ---
# The Status columns were sized for the long "Handed back: in sync with
# en-US" string; now that the text is the much shorter "Ready for
# handoff", the columns are narrowed to fit (an autofit-style resize).
# ColumnWidth is stored in the workbook on a pixel grid, so we assign the
# value that lands on the grid point closest to the narrowed width.
$ws_overview.Range("E1").ColumnWidth = 16.3333333333333
$ws_overview.Range("F1").ColumnWidth = 16.3333333333333
$ws_zhcn.Range("C1").ColumnWidth = 16.3333333333333
$ws_dede.Range("C1").ColumnWidth = 16.3333333333333
